$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.654.79'
$ws.Range("E2").Value = '  +2.50%  '
$ws.Range("D3").Value = '3.431.75'
$ws.Range("E3").Value = '  +3.00%  '
$ws.Range("E4").Value = '  -0.02%  '
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '407.04'
$ws.Range("D5").Style = $style

$ws.Range("E5").Value = '  +1.40%  '
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '130.33'
$ws.Range("D6").Style = $style

$ws.Range("E6").Value = '  +3.21%  '
$style = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.597'
$ws.Range("D7").Style = $style

$ws.Range("E7").Value = '  +1.06%  '
$ws.Range("E8").Value = '  -0.09%  '
$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.693'
$ws.Range("D9").Style = $style

$ws.Range("E9").Value = '  +5.26%  '
$ws.Range("E10").Value = '  +16.77%  '
$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.02'
$ws.Range("D11").Style = $style

$ws.Range("E11").Value = '  +2.68%  '
$ws.Range("E12").Value = '  +0.39%  '
$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.46'
$ws.Range("D13").Style = $style

$ws.Range("E13").Value = '  +2.12%  '
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.84'
$ws.Range("D14").Style = $style

$ws.Range("E14").Value = '  +2.93%  '
$ws.Range("D15").Value = '3.431.94'
$ws.Range("E15").Value = '  +3.17%  '
$ws.Range("D16").Value = '62.665.33'
$ws.Range("E16").Value = '  +2.67%  '
$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '11.53'
$ws.Range("D17").Style = $style

$ws.Range("E17").Value = '  +2.77%  '
$ws.Range("E18").Value = '  +1.22%  '
$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000163'
$ws.Range("D19").Style = $style

$ws.Range("E19").Value = '  +28.09%  '
$ws.Range("E20").Value = '  -0.26%  '
$ws.Range("E21").Value = '  +5.51%  '
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '315.28'
$ws.Range("D22").Style = $style

$ws.Range("E22").Value = '  +5.47%  '
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.86'
$ws.Range("D23").Style = $style

$ws.Range("E23").Value = '  +0.48%  '
$ws.Range("E24").Value = '  +1.92%  '
$style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.75'
$ws.Range("D25").Style = $style

$ws.Range("E25").Value = '  +0.38%  '
$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '29.77'
$ws.Range("D26").Style = $style

$ws.Range("E26").Value = '  +2.97%  '
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.18'
$ws.Range("D27").Style = $style

$ws.Range("E27").Value = '  -1.05%  '
$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.83'
$ws.Range("D28").Style = $style

$ws.Range("E28").Value = '  +5.98%  '
$ws.Range("E29").Value = '  +9.96%  '
$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '44.41'
$ws.Range("D30").Style = $style

$ws.Range("E30").Value = '  +8.33%  '
$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.173'
$ws.Range("D31").Style = $style

$ws.Range("E31").Value = '  +0.75%  '
$ws.Range("E32").Value = '  +1.32%  '
$ws.Range("E33").Value = '  +0.33%  '
$ws.Range("E34").Value = '  -0.05%  '
$ws.Range("E35").Value = '  +1.66%  '
$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '51.86'
$ws.Range("D36").Style = $style

$ws.Range("E36").Value = '  -0.38%  '
$ws.Range("E37").Value = '  +0.30%  '
$ws.Range("E38").Value = '  +2.11%  '
$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.322'
$ws.Range("D39").Style = $style

$ws.Range("E39").Value = '  +15.34%  '
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.33'
$ws.Range("D40").Style = $style

$ws.Range("E40").Value = '  -1.21%  '
$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '143.43'
$ws.Range("D41").Style = $style

$ws.Range("E41").Value = '  +5.21%  '
$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.126'
$ws.Range("D42").Style = $style

$ws.Range("E42").Value = '  +2.76%  '
$ws.Range("E43").Value = '  +0.66%  '
$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.93'
$ws.Range("D44").Style = $style

$ws.Range("E44").Value = '  +1.67%  '
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.91'
$ws.Range("D45").Style = $style

$ws.Range("E45").Value = '  +1.12%  '
$ws.Range("E46").Value = '  +0.25%  '
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '21.39'
$ws.Range("D47").Style = $style

$ws.Range("E47").Value = '  +0.37%  '
$ws.Range("D48").Value = '2.105.06'
$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.99'
$ws.Range("D49").Style = $style

$ws.Range("E49").Value = '  +7.44%  '
$ws.Range("E50").Value = '  -1.35%  '
$style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.09'
$ws.Range("D51").Style = $style

$ws.Range("E51").Value = '  +29.34%  '
